# Add newly documented wine regions to two worksheets.
$wb = $excel.ActiveWorkbook

# --- "Wine & food pairing" sheet: add Bordeaux (Red) pairing info ---
$wsPairing = $wb.Worksheets.Item("Wine & food pairing")
$wsPairing.Range("A11").Value = "Bordeaux (Red)"
$wsPairing.Range("B11").Value = "Lamb, game meat, beef, roasted, grilled, stewed; when tanin are round it goes well with chocolate desserts, brownies"
$wsPairing.Range("A12").Select()

# --- "European design. & varieties" sheet: add the two Mâcon regions ---
$wsRegions = $wb.Worksheets.Item("European design. & varieties")
$wsRegions.Select()
$wsRegions.Range("A39").Value = "Mâcon (Red), Bourgogne, France"
$wsRegions.Range("A40").Value = "Mâcon (White), Bourgogne, France"
$wsRegions.Range("B39").Value = "Gamay, Pinot noir"
$wsRegions.Range("B40").Value = "Chardonnay"
$wsRegions.Range("A41").Select()
